$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 22 new rows before the old "OpenMP" block (old row 34), pushing
# everything below down by 22 rows (old 34..53 -> new 56..75).
$ws.Rows("34:55").Insert()

# ---- Block 1: "Speedup στα 1000 Generations" (new rows 34-41) ----
$ws.Range("A34").Value = "Speedup στα 1000 Generations"

$ws.Range("A35").Value = "Gens=1000"
$ws.Range("B35").Value = "NPROB"
$ws.Range("C35").Value = 240
$ws.Range("D35").Value = 600
$ws.Range("E35").Value = 960
$ws.Range("F35").Value = 1440
$ws.Range("G35").Value = 2880
$ws.Range("H35").Value = 3600
$ws.Rows(35).RowHeight = 19.5

$ws.Range("A36").Value = "Processses"

$ws.Range("A37").Value = 1
$ws.Range("C37:H37").Value = 1
$ws.Range("C37:H37").NumberFormat = "0.000000"

$ws.Range("A38").Value = 4
$ws.Range("C38").Formula = "=C26/C27"
$ws.Range("D38").Formula = "=D26/D27"
$ws.Range("E38").Formula = "=E26/E27"
$ws.Range("F38").Formula = "=F26/F27"
$ws.Range("G38").Formula = "=G26/G27"
$ws.Range("H38").Formula = "=H26/H27"
$ws.Range("C38:H38").NumberFormat = "0.000000"

$ws.Range("A39").Value = 9
$ws.Range("C39").Formula = "=C26/C28"
$ws.Range("D39").Formula = "=D26/D28"
$ws.Range("E39").Formula = "=E26/E28"
$ws.Range("F39").Formula = "=F26/F28"
$ws.Range("G39").Formula = "=G26/G28"
$ws.Range("H39").Formula = "=H26/H28"
$ws.Range("C39:H39").NumberFormat = "0.000000"

$ws.Range("A40").Value = 16
$ws.Range("C40").Formula = "=C26/C29"
$ws.Range("D40").Formula = "=D26/D29"
$ws.Range("E40").Formula = "=E26/E29"
$ws.Range("F40").Formula = "=F26/F29"
$ws.Range("G40").Formula = "=G26/G29"
$ws.Range("H40").Formula = "=H26/H29"
$ws.Range("C40:H40").NumberFormat = "0.000000"

$ws.Range("A41").Value = 25
$ws.Range("C41").Formula = "=C26/C30"
$ws.Range("D41").Formula = "=D26/D30"
$ws.Range("E41").Formula = "=E26/E30"
$ws.Range("F41").Formula = "=F26/F30"
$ws.Range("G41").Formula = "=G26/G30"
$ws.Range("H41").Formula = "=H26/H30"
$ws.Range("C41:H41").NumberFormat = "0.000000"

# ---- Block 2: "Efficiency στα 1000 Generations" (new rows 44-51) ----
$ws.Range("A44").Value = "Efficiency στα 1000 Generations"
$ws.Range("D44").Value = "E=Tspeedup/Nprocesses"
$ws.Range("D34").Value = "S=Tserial/Tparal"

$ws.Range("A45").Value = "Gens=1000"
$ws.Range("B45").Value = "NPROB"
$ws.Range("C45").Value = 240
$ws.Range("D45").Value = 600
$ws.Range("E45").Value = 960
$ws.Range("F45").Value = 1440
$ws.Range("G45").Value = 2880
$ws.Range("H45").Value = 3600

$ws.Range("A46").Value = "Processses"

$ws.Range("A47").Value = 1
$ws.Range("C47:H47").Value = 1
$ws.Range("C47:H47").NumberFormat = "0.000000"

$ws.Range("A48").Value = 4
$ws.Range("C48").Formula = "=C38/A48"
$ws.Range("D48").Formula = "=D38/A48"
$ws.Range("E48").Formula = "=E38/A48"
$ws.Range("F48").Formula = "=F38/A48"
$ws.Range("G48").Formula = "=G38/A48"
$ws.Range("H48").Formula = "=H38/A48"
$ws.Range("C48:H48").NumberFormat = "0.000000"

$ws.Range("A49").Value = 9
$ws.Range("C49").Formula = "=C39/A49"
$ws.Range("D49").Formula = "=D39/A49"
$ws.Range("E49").Formula = "=E39/A49"
$ws.Range("F49").Formula = "=F39/A49"
$ws.Range("G49").Formula = "=G39/A49"
$ws.Range("H49").Formula = "=H39/A49"
$ws.Range("C49:H49").NumberFormat = "0.000000"

$ws.Range("A50").Value = 16
$ws.Range("C50").Formula = "=C40/A50"
$ws.Range("D50").Formula = "=D40/A50"
$ws.Range("E50").Formula = "=E40/A50"
$ws.Range("F50").Formula = "=F40/A50"
$ws.Range("G50").Formula = "=G40/A50"
$ws.Range("H50").Formula = "=H40/A50"
$ws.Range("C50:H50").NumberFormat = "0.000000"

$ws.Range("A51").Value = 25
$ws.Range("C51").Formula = "=C41/A51"
$ws.Range("D51").Formula = "=D41/A51"
$ws.Range("E51").Formula = "=E41/A51"
$ws.Range("F51").Formula = "=F41/A51"
$ws.Range("G51").Formula = "=G41/A51"
$ws.Range("H51").Formula = "=H41/A51"
$ws.Range("C51:H51").NumberFormat = "0.000000"

# The rows that got pushed down (old row 35, now row 57) lose the custom
# row height that used to belong to the "OpenMP" header row; only the new
# "Speedup" header row (35) keeps a custom height now.
$ws.Rows(57).AutoFit()

# ---- View state: selection moves to F53 (scroll position best-effort) ----
[void]$ws.Range("F53").Select()
